$d = $word.ActiveDocument

# Make sure edits are applied as plain content changes, not tracked revisions.
$d.TrackRevisions = $false

# ------------------------------------------------------------------
# 1) HPCDATAMGM-1556 bullet: body text is replaced with the text that
#    used to belong to the HPCDATAMGM-1560 bullet (reworded slightly:
#    "Browse screen" -> "Browse page", "display archive" -> "displayed
#    archive").
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Enhanced the Browse screen of the DME web application to show the data hierarchy of the display archive (base path) through an info icon. Previously, the user had to make a REST API call to determine the data hierarchy for a specific base path.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Enhanced the Browse page of the DME web application to show the data hierarchy of the displayed archive (base path) through an info icon. Previously, users had to navigate to the Register Collection page to obtain this information. ",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) HPCDATAMGM-1560 bullet: body text is replaced with new wording
#    about mandatory metadata associated with each collection type.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Enhanced the Browse screen of the DME web application to show the mandatory metadata at each level of the data hierarchy for the displayed archive through an info icon. Previously, users had to navigate to the Register Collection page to obtain this information. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Enhanced the Browse page of the DME web application to show the mandatory metadata associated with each collection type in the hierarchy of the displayed archive (base path). Previously, the user had to navigate to the Register Collection page and select each item from the Collection Type dropdown to obtain this information.",
    2) | Out-Null

# Add the leading space that now precedes "HPCDATAMGM-1560" in the
# revised document.
$d.Content.Find.Execute(
    "HPCDATAMGM-1560:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " HPCDATAMGM-1560:",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) HPCDATAMGM-1567 bullet: reword "is now applicable to non-admins.
#    only." -> "now applies only to non-administrators."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " is now applicable to non-admins. only. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " now applies only to non-administrators. ",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) HPCDATAMGM-1555 bullet: tidy the spacing inside the parentheses.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "( accessed from Reports tab )",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "(accessed from Reports tab)",
    2) | Out-Null

# ------------------------------------------------------------------
# 5) HPCDATAMGM-1542 bullet: wrap the issue number in a bookmark
#    (mirrors a Word-generated "_Hlk" bookmark from the review pass).
# ------------------------------------------------------------------
$bmRange = $d.Content
$bmRange.Find.Execute("HPCDATAMGM-1542", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_Hlk96681428", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 6) HPCDATAMGM-1525 bullet: add a comma after "Previously".
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Previously only the failed status was provided.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Previously, only the failed status was provided.",
    2) | Out-Null

# ------------------------------------------------------------------
# 7) HPCDATAMGM-1551 bullet: several small wording tweaks.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Added DME managed thread pool to limit the number of threads setup during streaming transfers (i.e. downloads or uploads to Cloudian or AWS S3). This replaces the default AWS thread pool, providing the ability to control the number of parts",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Added DME managed thread pool to limit the number of threads during streaming transfers (such as downloads or uploads to Cloudian or AWS S3). This replaces the default AWS thread pool, providing the ability to control the number of parts",
    2) | Out-Null
